# Add two new worksheets ("same_elements" and "partly_same") at the end of
# the workbook, populate them with the "same elements" / "partly same
# elements" sorting-benchmark results, and leave the new last-added sheet
# ("same_elements") selected/active, matching the author's commit
# "same elements array researches has been added".

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$s5 = $wb.Worksheets.Add($null, $lastSheet)
$s5.Name = "same_elements"

$s6 = $wb.Worksheets.Add($null, $s5)
$s6.Name = "partly_same"

# --- same_elements ---------------------------------------------------
$s5.Range("B1").Value = 5
$s5.Range("C1").Value = 50
$s5.Range("D1").Value = 500
$s5.Range("E1").Value = 5000
$s5.Range("F1").Value = 50000
$s5.Range("G1").Value = 500000

$s5.Range("A2").Value = "byte"
$s5.Range("B2").Value = 0
$s5.Range("C2").Value = 0
$s5.Range("D2").Value = 0.042544
$s5.Range("E2").Value = 4.577781
$s5.Range("F2").Value = 0
$s5.Range("G2").Value = 0

$s5.Range("A3").Value = "int"
$s5.Range("B3").Value = 0
$s5.Range("C3").Value = 0
$s5.Range("D3").Value = 0.041543
$s5.Range("E3").Value = 4.132817
$s5.Range("F3").Value = 0
$s5.Range("G3").Value = 0

$s5.Range("A4").Value = "string"
$s5.Range("B4").Value = 0
$s5.Range("C4").Value = 0
$s5.Range("D4").Value = 0.050053
$s5.Range("E4").Value = 4.099128
$s5.Range("F4").Value = 0
$s5.Range("G4").Value = 0

$s5.Range("A5").Value = "date"
$s5.Range("B5").Value = 0
$s5.Range("C5").Value = 0.000501
$s5.Range("D5").Value = 0.051053
$s5.Range("E5").Value = 4.162854
$s5.Range("F5").Value = 0
$s5.Range("G5").Value = 0

# --- partly_same -------------------------------------------------------
$s6.Range("B1").Value = 5
$s6.Range("C1").Value = 50
$s6.Range("D1").Value = 500
$s6.Range("E1").Value = 5000
$s6.Range("F1").Value = 50000
$s6.Range("G1").Value = 500000

$s6.Range("A2").Value = "byte"
$s6.Range("B2").Value = 0
$s6.Range("C2").Value = 0
$s6.Range("D2").Value = 0.039541
$s6.Range("E2").Value = 3.976907
$s6.Range("F2").Value = 0
$s6.Range("G2").Value = 0

$s6.Range("A3").Value = "int"
$s6.Range("B3").Value = 0
$s6.Range("C3").Value = 0.000502
$s6.Range("D3").Value = 0.04855
$s6.Range("E3").Value = 4.081778
$s6.Range("F3").Value = 0
$s6.Range("G3").Value = 0

$s6.Range("A4").Value = "string"
$s6.Range("B4").Value = 0
$s6.Range("C4").Value = 0.0005
$s6.Range("D4").Value = 0.040649
$s6.Range("E4").Value = 4.015055
$s6.Range("F4").Value = 0
$s6.Range("G4").Value = 0

$s6.Range("A5").Value = "date"
$s6.Range("B5").Value = 0
$s6.Range("C5").Value = 0
$s6.Range("D5").Value = 0.035536
$s6.Range("E5").Value = 4.077259
$s6.Range("F5").Value = 0
$s6.Range("G5").Value = 0

# --- leave "same_elements" selected/active, matching the commit --------
$s5.Select()
$s5.Range("H13").Select()
